$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Re-map the header row (A1:K1) onto its new layout.
#    Using Range.Copy (instead of .Value = "...") preserves the
#    original shared-string / text typing so numeric-looking labels
#    like "14" do not get silently reinterpreted as numbers, and so
#    no new cell styles get minted.
#    Because several header cells trade places with each other, the
#    old row is staged in a scratch row (row 20) first so that every
#    copy reads an untouched source cell.
# ------------------------------------------------------------------
$ws.Range("A1:K1").Copy($ws.Range("A20"))

$ws.Range("G20").Copy($ws.Range("B1"))   # HOAc
$ws.Range("B20").Copy($ws.Range("C1"))   # 14
$ws.Range("I20").Copy($ws.Range("D1"))   # H2
$ws.Range("E20").Copy($ws.Range("E1"))   # Formaldehyde
$ws.Range("A20").Copy($ws.Range("F1"))   # L
$ws.Range("H20").Copy($ws.Range("G1"))   # CO
$ws.Range("J20").Copy($ws.Range("H1"))   # Water
$ws.Range("D20").Copy($ws.Range("I1"))   # 31-ol
$ws.Range("K20").Copy($ws.Range("J1"))   # CO2

# A1 no longer carries a header label.
$ws.Range("A1").Clear()

# Drop the scratch staging row entirely.
$ws.Rows("20:20").Delete()

# ------------------------------------------------------------------
# 2) Update the data row (row 2) values.
# ------------------------------------------------------------------
$ws.Range("D2").Value = 1607318.565012615
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 60978.88933325125
$ws.Range("H2").Value = -1
$ws.Range("I2").Value = -1
$ws.Range("J2").Value = 282318.2582178094

# ------------------------------------------------------------------
# 3) Remove the now-obsolete rows 3:5 and column K.
# ------------------------------------------------------------------
$ws.Rows("3:5").Delete()
$ws.Columns("K").Delete()
